# Fixing "beating the red light" issue:
#  1) Add a space before "&A" in the odd header and before "Page" in the
#     odd footer for the Impediment_Rule and Traffic_Light_Rule sheets.
#  2) Update the fuzzy membership function boundary values on the
#     fuzzy_values sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Header / footer tweaks -------------------------------------------
foreach ($sheetName in @("Impediment_Rule", "Traffic_Light_Rule")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12 &A'
    $ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12 Page &P'
}

# --- 2) fuzzy_values data updates ------------------------------------------
$ws = $wb.Worksheets.Item("fuzzy_values")

$ws.Range("D2").Value = 15.5699101643917
$ws.Range("E2").Value = 145.1814757813382

$ws.Range("B3").Value = 15.5699101643917
$ws.Range("C3").Value = 145.1814757813382
$ws.Range("D3").Value = 149.9585583888527

$ws.Range("B4").Value = 145.1814757813382
$ws.Range("C4").Value = 149.9585583888527

$ws.Range("D5").Value = 17.16511041311179
$ws.Range("E5").Value = 114.6172243513307

$ws.Range("B6").Value = 17.16511041311179
$ws.Range("C6").Value = 114.6172243513307
$ws.Range("D6").Value = 201.6501654319696

$ws.Range("B7").Value = 114.6172243513307
$ws.Range("C7").Value = 201.6501654319696

$ws.Range("B8").Value = 1.499543122584523
$ws.Range("C8").Value = 20.37558392861101

$ws.Range("D9").Value = 1.499543122584523
$ws.Range("E9").Value = 20.37558392861101

$ws.Range("B10").Value = 10.56784383741841
$ws.Range("C10").Value = 18.40942901035632

$ws.Range("D11").Value = 10.56784383741841
$ws.Range("E11").Value = 18.40942901035632

$ws.Range("C13").Value = 0.2961582525032522
$ws.Range("D13").Value = 0.875326724110622

$ws.Range("B14").Value = 0.2961582525032522
$ws.Range("C14").Value = 0.875326724110622
$ws.Range("D14").Value = 1.851474461442805

$ws.Range("B15").Value = 0.875326724110622
$ws.Range("C15").Value = 1.851474461442805
